$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "265.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.65%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.705"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.06%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06079"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.30%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8501"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9052"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.22%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1406"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.18%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.73%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07101"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03141"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.45%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09028"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006056"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.62%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.97%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.458"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.06%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.170"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3091"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.36%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1300"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.65%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.106"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.19%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04240"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.57%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001182"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.87%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004136"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.77%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.07%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03924"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1115"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.24%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004171"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.25%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.33%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01154"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-16.58%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005103"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.63%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1245"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-25.31%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
